$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the last row of data (based on column A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "r585"
$ws.Cells.Item($newRow, 2).Value = "cameron"
$ws.Cells.Item($newRow, 3).Value = "this is with a feedback form"
$ws.Cells.Item($newRow, 4).Value = "2025-10-01 16:26:29"
